$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Book2Children"
$ws.Range("B3").Value = "Delhi"
$ws.Range("C3").Value = "Chennai"
$ws.Range("F3").Value = "'2"

$ws.Range("C4").Select()
